$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '42.827.62'
$ws.Range("E2").Value = '  -4.83%  '

# Row 3
$ws.Range("D3").Value = '2.220.56'
$ws.Range("E3").Value = '  -6.08%  '

# Row 4
$ws.Range("E4").Value = '  -0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '316.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.03'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.00%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.595'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.58%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.24%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.564'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.24%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.06%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.89'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.87%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0828'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -9.22%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.83'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.74%  '

# Row 14
$ws.Range("E14").Value = '  -3.16%  '

# Row 15
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.861'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -11.02%  '

# Row 16
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.557.97'
$ws.Range("E16").Value = '  -6.09%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.28'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.85%  '

# Row 18
$ws.Range("D18").Value = '2.219.40'
$ws.Range("E18").Value = '  -6.27%  '

# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '42.740.74'
$ws.Range("E19").Value = '  -4.94%  '

# Row 20
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.10%  '

# Row 21
$ws.Range("E21").Value = '  -8.47%  '

# Row 22
$ws.Range("E22").Value = '  -10.54%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.39%  '

# Row 24
$ws.Range("E24").Value = '  -8.55%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '236.65'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.54%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.14'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.86%  '

# Row 27
$ws.Range("E27").Value = '  -0.13%  '

# Row 28
$ws.Range("E28").Value = '  -8.51%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -10.65%  '

# Row 31
$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0889'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.24%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.49'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.04%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -7.42%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '156.66'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.85%  '

# Row 35
$ws.Range("E35").Value = '  -6.43%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.17'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.36%  '

# Row 37
$ws.Range("E37").Value = '  +14.12%  '

# Row 38
$ws.Range("E38").Value = '  -5.28%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.66%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.53%  '

# Row 41
$ws.Range("E41").Value = '  -10.65%  '

# Row 42
$ws.Range("E42").Value = '  -7.41%  '

# Row 43
$ws.Range("D43").Value = '1.934.82'
$ws.Range("E43").Value = '  +3.68%  '

# Row 44
$ws.Range("E44").Value = '  -0.12%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.23'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.71%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.54'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.46%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.208'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.52%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.03%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '60.77'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -12.29%  '

# Row 50
$ws.Range("E50").Value = '  -6.77%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.871'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +19.80%  '
